$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Menus_and_forms")

# Update header row E1:H1 - replace old filepath-based column names with new names
$ws.Range("E1").Value = "image_en"
$ws.Range("F1").Value = "audio_en"
$ws.Range("G1").Value = "image_fra"
$ws.Range("H1").Value = "audio_fra"

# Update the saved selection on this sheet to E2
$ws.Range("E2").Select()
